$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-formatted cells to stay as text (avoid Excel auto-numeric conversion)
$textCells = @("D5", "D6", "D10", "D12", "D14", "D17", "D19", "D22", "D24", "D26", "D27", "D28", "D30", "D31", "D33", "D35", "D36", "D37", "D44", "D45", "D46", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply all cell value changes
$ws.Cells.Item(2, 4).Value = '42.670.47'
$ws.Cells.Item(2, 5).Value = '  -0.60%  '
$ws.Cells.Item(3, 4).Value = '2.294.26'
$ws.Cells.Item(3, 5).Value = '  -0.38%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$ws.Cells.Item(5, 4).Value = '302.11'
$ws.Cells.Item(5, 5).Value = '  +0.67%  '
$ws.Cells.Item(6, 4).Value = '96.03'
$ws.Cells.Item(6, 5).Value = '  -1.77%  '
$ws.Cells.Item(7, 5).Value = '  -1.95%  '
$ws.Cells.Item(8, 5).Value = '  +0.09%  '
$ws.Cells.Item(9, 5).Value = '  -2.27%  '
$ws.Cells.Item(10, 4).Value = '34.74'
$ws.Cells.Item(10, 5).Value = '  -3.12%  '
$ws.Cells.Item(11, 5).Value = '  -0.82%  '
$ws.Cells.Item(12, 4).Value = '18.55'
$ws.Cells.Item(12, 5).Value = '  +4.69%  '
$ws.Cells.Item(13, 5).Value = '  +2.40%  '
$ws.Cells.Item(14, 4).Value = '6.84'
$ws.Cells.Item(14, 5).Value = '  +0.81%  '
$ws.Cells.Item(15, 4).Value = '2.651.88'
$ws.Cells.Item(15, 5).Value = '  -0.22%  '
$ws.Cells.Item(16, 4).Value = '2.310.64'
$ws.Cells.Item(16, 5).Value = '  -0.10%  '
$ws.Cells.Item(17, 4).Value = '0.772'
$ws.Cells.Item(17, 5).Value = '  -0.85%  '
$ws.Cells.Item(18, 4).Value = '42.601.76'
$ws.Cells.Item(18, 5).Value = '  -0.66%  '
$ws.Cells.Item(19, 4).Value = '12.83'
$ws.Cells.Item(19, 5).Value = '  +1.79%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0891'
$ws.Cells.Item(20, 5).Value = '  -1.77%  '
$ws.Cells.Item(21, 5).Value = '  -1.84%  '
$ws.Cells.Item(22, 4).Value = '67.04'
$ws.Cells.Item(22, 5).Value = '  -1.36%  '
$ws.Cells.Item(23, 5).Value = '  -2.40%  '
$ws.Cells.Item(24, 4).Value = '2.12'
$ws.Cells.Item(24, 5).Value = '  -0.72%  '
$ws.Cells.Item(25, 5).Value = '  +0.08%  '
$ws.Cells.Item(26, 4).Value = '2.39'
$ws.Cells.Item(27, 4).Value = '24.54'
$ws.Cells.Item(27, 5).Value = '  -2.14%  '
$ws.Cells.Item(28, 4).Value = '167.50'
$ws.Cells.Item(28, 5).Value = '  +0.76%  '
$ws.Cells.Item(29, 5).Value = '  +0.42%  '
$ws.Cells.Item(30, 4).Value = '8.97'
$ws.Cells.Item(30, 5).Value = '  -0.71%  '
$ws.Cells.Item(31, 4).Value = '32.79'
$ws.Cells.Item(31, 5).Value = '  -0.06%  '
$ws.Cells.Item(32, 5).Value = '  +0.06%  '
$ws.Cells.Item(33, 4).Value = '17.73'
$ws.Cells.Item(33, 5).Value = '  +1.34%  '
$ws.Cells.Item(34, 5).Value = '  -1.19%  '
$ws.Cells.Item(35, 4).Value = '4.44'
$ws.Cells.Item(35, 5).Value = '  -6.38%  '
$ws.Cells.Item(36, 4).Value = '2.35'
$ws.Cells.Item(36, 5).Value = '  -1.84%  '
$ws.Cells.Item(37, 4).Value = '0.0685'
$ws.Cells.Item(37, 5).Value = '  -0.20%  '
$ws.Cells.Item(38, 5).Value = '  -0.52%  '
$ws.Cells.Item(39, 5).Value = '  -1.76%  '
$ws.Cells.Item(41, 5).Value = '  -2.86%  '
$ws.Cells.Item(42, 4).Value = '1.992.56'
$ws.Cells.Item(42, 5).Value = '  -0.40%  '
$ws.Cells.Item(43, 5).Value = '  -2.04%  '
$ws.Cells.Item(44, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(44, 4).Value = '2.14'
$ws.Cells.Item(44, 5).Value = '  -0.24%  '
$ws.Cells.Item(45, 2).Value = 'FraxShare'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(45, 4).Value = '10.13'
$ws.Cells.Item(45, 5).Value = '  -0.07%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Value = '18.21'
$ws.Cells.Item(46, 5).Value = '  +5.64%  '
$ws.Cells.Item(47, 5).Value = '  -0.70%  '
$ws.Cells.Item(48, 2).Value = 'MultiversX'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Cells.Item(48, 4).Value = '53.48'
$ws.Cells.Item(48, 5).Value = '  +0.22%  '
$ws.Cells.Item(49, 2).Value = 'HuobiToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(49, 4).Value = '2.83'
$ws.Cells.Item(49, 5).Value = '  +4.19%  '
$ws.Cells.Item(50, 4).Value = '2.518.74'
$ws.Cells.Item(50, 5).Value = '  -0.22%  '
$ws.Cells.Item(51, 4).Value = '70.73'
$ws.Cells.Item(51, 5).Value = '  -2.18%  '

# Restore default style on the text-forced cells (removes the temporary text format)
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
